$d = $word.ActiveDocument

# 1. "Indie Panto Pop" -> "Indie Panto Party Pop"
$d.Content.Find.Execute("Indie Panto Pop", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Indie Panto Party Pop", 2) | Out-Null

# 2. Add a new closing paragraph after the "2019 saw Duck Thieves..." paragraph,
#    before the bookmark end / trailing empty paragraph.
$paras = $d.Paragraphs
$lastContentPara = $paras.Item(7)
$lastContentPara.Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item(8).Range.Text = "Duck Thieves are looking to release a new EP in 2023"
